# Apply the "doc: update tutorial model files" edit to the RAVEN tutorial
# workbook:
#   * RXNS sheet: a new "CONFIDENCE SCORE" style #NUM! error value is
#     written for every reaction data row, and the two section-separator
#     rows (20 and 26) get their usual blank/formatted cells.
#   * MODEL sheet: header row is rewritten (DESCRIPTION -> NAME, with a new
#     TAXONOMY column inserted right after NAME), and the sample data row is
#     trimmed down to just the id / name / notes values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# RXNS sheet (sheet1)
# ---------------------------------------------------------------------
$rxns = $wb.Worksheets.Item("RXNS")

$dataRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,21,22,23,24,25,27)
foreach ($r in $dataRows) {
    $cell = $rxns.Range("P$r")
    $cell.NumberFormat = $rxns.Range("D$r").NumberFormat
    $cell.Value = "#NUM!"
}

$blankRows = @(20,26)
$blankCols = @("B","C","D","J","P")
foreach ($r in $blankRows) {
    foreach ($col in $blankCols) {
        $cell = $rxns.Range("$col$r")
        $cell.NumberFormat = $rxns.Range("${col}21").NumberFormat
        $cell.Value = ""
    }
}

# ---------------------------------------------------------------------
# MODEL sheet (sheet4)
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("MODEL")

# Header row: DESCRIPTION -> NAME, insert TAXONOMY, shift the remaining
# headers one column to the right (NOTES stays put at the end).
$model.Range("C1").Value = "NAME"
$model.Range("D1").Value = "TAXONOMY"
$model.Range("E1").Value = "DEFAULT LOWER"
$model.Range("F1").Value = "DEFAULT UPPER"
$model.Range("G1").Value = "CONTACT GIVEN NAME"
$model.Range("H1").Value = "CONTACT FAMILY NAME"
$model.Range("I1").Value = "CONTACT EMAIL"
$model.Range("J1").Value = "ORGANIZATION"
$model.Range("K1").Value = "NOTES"

# Data row: keep id/name/notes, clear the rest of the sample values.
$model.Range("D2").Value = ""
$model.Range("E2").Value = ""
$model.Range("F2").Value = ""
$model.Range("G2").Value = ""
$model.Range("H2").Value = ""
$model.Range("I2").Value = ""
$model.Range("J2").Value = ""
